$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (45203) on every data row
# (rows 2-398). Update it to the new date serial value 45205.
$lastRow = 398
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45205
